# Apply Search Module locator updates to the "Locators" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locators")

# Rows 27-29: switch LocatorType from css id-selector to xpath id-selector
$ws.Cells.Item(27, 4).Value = "xpath"
$ws.Cells.Item(27, 5).Value = "//input[@id='firstname-input']"

$ws.Cells.Item(28, 4).Value = "xpath"
$ws.Cells.Item(28, 5).Value = "//input[@id='lastname-input']"

$ws.Cells.Item(29, 4).Value = "xpath"
$ws.Cells.Item(29, 5).Value = "//input[@id='personid-input']"

# Rows 31-32: reuse the same selectize-input xpath (index [1]) as row 30
$ws.Cells.Item(31, 5).Value = "(//div[@class='selectize-input items not-full has-options']/input)[1]"
$ws.Cells.Item(32, 5).Value = "(//div[@class='selectize-input items not-full has-options']/input)[1]"

# Row 33: new locator for the first search-result row's name link
$ws.Cells.Item(33, 1).Value = "Tgl"
$ws.Cells.Item(33, 2).Value = "SearchPage"
$ws.Cells.Item(33, 3).Value = "Tgl_firstrow_name"
$ws.Cells.Item(33, 4).Value = "xpath"
$ws.Cells.Item(33, 5).Value = "//tbody[@data-hook='results']/tr[1]/td/a"

# Update the sheet's active selection to C33
$ws.Activate()
$ws.Range("C33").Select()
